# Generate Report for Handback
# Updates the timestamp values recorded on the handback-status report:
#  - Overview sheet: "Latest HO Xliff Generate Date" for the first file
#  - zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first file
#  - de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first file

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-23 23:01:21"

# zh-cn!H2 - Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-23 23:01:15"
# zh-cn!K2 - Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-23 23:01:33"

# de-de!H2 - Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-08-23 23:01:21"
# de-de!K2 - Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-23 23:01:40"
